# Atualização de bases das ligas, do dia: 2024-02-02 às 20:58
#
# Two pairs of match rows had their betting-odds data swapped back to the
# correct rows (the "id" / teams / odds columns move, while the leading
# row-index column A stays put):
#   - rows 233 <-> 234
#   - rows 236 <-> 239   and   rows 237 <-> 238

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param(
        [int]$Row,
        [double]$B,
        [string]$F,
        [string]$G,
        [double]$H,
        [double]$I,
        [string]$J,
        [double]$K,
        [double]$L,
        [double]$M,
        [double]$N,
        [double]$O,
        [double]$P,
        [double]$Q,
        [double]$R,
        [double]$S,
        [double]$T,
        [double]$U,
        [double]$V,
        [double]$W,
        [double]$X,
        [double]$Y,
        [double]$Z,
        [double]$AA,
        [double]$AB,
        [double]$AC
    )

    $ws.Cells.Item($Row, 2).Value  = $B   # B  id
    $ws.Cells.Item($Row, 6).Value  = $F   # F  HomeTeam
    $ws.Cells.Item($Row, 7).Value  = $G   # G  AwayTeam
    $ws.Cells.Item($Row, 8).Value  = $H   # H  FTHG
    $ws.Cells.Item($Row, 9).Value  = $I   # I  FTAG
    $ws.Cells.Item($Row, 10).Value = $J   # J  FTR
    $ws.Cells.Item($Row, 11).Value = $K   # K  oddH_op
    $ws.Cells.Item($Row, 12).Value = $L   # L  oddD_op
    $ws.Cells.Item($Row, 13).Value = $M   # M  oddA_op
    $ws.Cells.Item($Row, 14).Value = $N   # N  oddH
    $ws.Cells.Item($Row, 15).Value = $O   # O  oddD
    $ws.Cells.Item($Row, 16).Value = $P   # P  oddA
    $ws.Cells.Item($Row, 17).Value = $Q   # Q  Ah
    $ws.Cells.Item($Row, 18).Value = $R   # R  oddAHH
    $ws.Cells.Item($Row, 19).Value = $S   # S  oddAHA
    $ws.Cells.Item($Row, 20).Value = $T   # T  AhOU
    $ws.Cells.Item($Row, 21).Value = $U   # U  oddAHOver
    $ws.Cells.Item($Row, 22).Value = $V   # V  oddAHUnder
    $ws.Cells.Item($Row, 23).Value = $W   # W  PLH
    $ws.Cells.Item($Row, 24).Value = $X   # X  PLD
    $ws.Cells.Item($Row, 25).Value = $Y   # Y  PLA
    $ws.Cells.Item($Row, 26).Value = $Z   # Z  PL_Ahh
    $ws.Cells.Item($Row, 27).Value = $AA  # AA PL_Aha
    $ws.Cells.Item($Row, 28).Value = $AB  # AB PL_AhOver
    $ws.Cells.Item($Row, 29).Value = $AC  # AC PL_AhUnder
}

# NOTE: parameters are passed positionally (not with -Name value) because
# negative numeric values break named-parameter binding in this shell.

# Row 233 (index 231) <- old row 234's data
Set-Row 233 7559469 "Montevideo Wanderers" "Penarol" 0 0 "D" 4.75 3.4 1.7 2.7 3.2 2.45 0 2.05 1.8 2.5 1.975 1.875 -1 2.2 -1 0 0 -1 0.875

# Row 234 (index 232) <- old row 233's data
Set-Row 234 7559468 "Liverpool Montevideo" "CA River Plate" 2 1 "H" 1.7 3 5.75 1.833 3.2 4.5 -0.5 1.925 1.925 2.25 2.025 1.825 0.833 -1 -1 0.925 -1 1.025 -1

# Row 236 (index 234) <- old row 239's data
Set-Row 236 7013886 "Racing Club de Montevideo" "Cerro" 0 1 "A" 2.25 3.1 3.25 2.25 2.875 3.5 -0.25 1.95 1.9 2 1.925 1.925 -1 -1 2.5 -1 0.8999999999999999 -1 0.925

# Row 237 (index 235) <- old row 238's data
Set-Row 237 7013885 "La Luz" "Atletico Fenix Montevideo" 0 2 "A" 3 3 2.4 2.9 2.75 2.6 0 2.025 1.825 2 2.025 1.825 -1 -1 1.6 -1 0.825 0 0

# Row 238 (index 236) <- old row 237's data
Set-Row 238 7013409 "Nacional De Football" "Torque" 1 1 "D" 1.666 3.9 4.5 1.615 4 4.75 -0.75 1.8 2.05 2.75 1.95 1.9 -1 3 -1 -1 1.05 -1 0.8999999999999999

# Row 239 (index 237) <- old row 236's data
Set-Row 239 7013702 "Defensor Sporting" "Danubio" 0 2 "A" 1.8 3.6 4.2 1.8 3.6 4.2 -0.75 2.05 1.8 2.25 1.85 2 -1 -1 3.2 -1 0.8 -0.5 0.5
